$wb = $excel.ActiveWorkbook

$alc = $wb.Worksheets.Item("ALC")
$alc.Range("H15").Value = 475.2
$alc.Range("I15").Value = 475.2
$alc.Range("K15").Value = 1425.6
$alc.Range("M15").Value = -1256.6

$alc.Range("H80").Value = 1050.5
$alc.Range("I80").Value = 1075
$alc.Range("J80").Value = 1001.5
$alc.Range("K80").Value = 3225
$alc.Range("L80").Value = 3004.5
$alc.Range("M80").Value = -2227
$alc.Range("N80").Value = -5000.5

$alc.Range("H83").Value = 1050.5
$alc.Range("I83").Value = 1075
$alc.Range("J83").Value = 1001.5
$alc.Range("K83").Value = 9675
$alc.Range("L83").Value = 9013.5
$alc.Range("M83").Value = -4683
$alc.Range("N83").Value = -18997.5

$alc.Range("H103").Value = 838.73914
$alc.Range("I103").Value = 750
$alc.Range("J103").Value = 847.1905
$alc.Range("K103").Value = 2250
$alc.Range("L103").Value = 2541.5715
$alc.Range("M103").Value = -1664
$alc.Range("N103").Value = -3713.5715

$arm = $wb.Worksheets.Item("ARM")
$arm.Range("H41").Value = 4118
$arm.Range("I41").Value = 1897.5
$arm.Range("J41").Value = 13000
$arm.Range("K41").Value = 1897.5
$arm.Range("L41").Value = 13000
$arm.Range("M41").Value = -1483.5
$arm.Range("N41").Value = -13828

$arm.Range("H61").Value = 3188.0588
$arm.Range("I61").Value = 2243.7144
$arm.Range("J61").Value = 4713.5386
$arm.Range("K61").Value = 2243.7144
$arm.Range("L61").Value = 4713.5386
$arm.Range("M61").Value = -2031.7144
$arm.Range("N61").Value = -5137.5386

$arm.Range("H74").Value = 1571.5264
$arm.Range("I74").Value = 1250.7727
$arm.Range("J74").Value = 2012.5625
$arm.Range("K74").Value = 1250.7727
$arm.Range("L74").Value = 2012.5625
$arm.Range("M74").Value = -376.7727
$arm.Range("N74").Value = -3760.5625

$arm.Range("H77").Value = 1571.5264
$arm.Range("I77").Value = 1250.7727
$arm.Range("J77").Value = 2012.5625
$arm.Range("K77").Value = 6253.863499999999
$arm.Range("L77").Value = 10062.8125
$arm.Range("M77").Value = -1885.863499999999
$arm.Range("N77").Value = -18798.8125

$arm.Range("H88").Value = 2499.182
$arm.Range("I88").Value = 2347.75
$arm.Range("J88").Value = 2585.7144
$arm.Range("K88").Value = 2347.75
$arm.Range("L88").Value = 2585.7144
$arm.Range("M88").Value = -1941.75
$arm.Range("N88").Value = -3397.7144

$arm.Range("H91").Value = 2499.182
$arm.Range("I91").Value = 2347.75
$arm.Range("J91").Value = 2585.7144
$arm.Range("K91").Value = 2347.75
$arm.Range("L91").Value = 2585.7144
$arm.Range("M91").Value = -943.75
$arm.Range("N91").Value = -5393.7144

$arm.Range("H131").Value = 43453.848
$arm.Range("J131").Value = 43453.848
$arm.Range("L131").Value = 43453.848
$arm.Range("N131").Value = -53533.848

$arm.Range("H132").Value = 2780755.5
$arm.Range("I132").Value = 1708.88
$arm.Range("J132").Value = 9096770
$arm.Range("K132").Value = 5126.64
$arm.Range("L132").Value = 27290310
$arm.Range("M132").Value = -2596.64
$arm.Range("N132").Value = -27295370

$arm.Range("H136").Value = 3188.0588
$arm.Range("I136").Value = 2243.7144
$arm.Range("J136").Value = 4713.5386
$arm.Range("K136").Value = 6731.1432
$arm.Range("L136").Value = 14140.6158
$arm.Range("M136").Value = -4181.1432
$arm.Range("N136").Value = -19240.6158

$bsm = $wb.Worksheets.Item("BSM")
$bsm.Range("H86").Value = 2133.1333
$bsm.Range("I86").Value = 1967.6666
$bsm.Range("J86").Value = 2243.4443
$bsm.Range("K86").Value = 1967.6666
$bsm.Range("L86").Value = 2243.4443
$bsm.Range("M86").Value = -844.6666
$bsm.Range("N86").Value = -4489.4443

$bsm.Range("H89").Value = 2133.1333
$bsm.Range("I89").Value = 1967.6666
$bsm.Range("J89").Value = 2243.4443
$bsm.Range("K89").Value = 9838.333
$bsm.Range("L89").Value = 11217.2215
$bsm.Range("M89").Value = -4222.333000000001
$bsm.Range("N89").Value = -22449.2215

$bsm.Range("H107").Value = 1099.1
$bsm.Range("I107").Value = 984.8823
$bsm.Range("J107").Value = 1248.4615
$bsm.Range("K107").Value = 984.8823
$bsm.Range("L107").Value = 1248.4615
$bsm.Range("M107").Value = 935.1177
$bsm.Range("N107").Value = -5088.461499999999

$bsm.Range("H134").Value = 2138.5293
$bsm.Range("I134").Value = 1892.4
$bsm.Range("J134").Value = 2822.2222
$bsm.Range("K134").Value = 5677.200000000001
$bsm.Range("L134").Value = 8466.6666
$bsm.Range("M134").Value = -3142.200000000001
$bsm.Range("N134").Value = -13536.6666

$cul = $wb.Worksheets.Item("CUL")
$cul.Range("H113").Value = 1304941.4
$cul.Range("I113").Value = 1389499.8
$cul.Range("J113").Value = 1000531.4
$cul.Range("K113").Value = 4168499.4
$cul.Range("L113").Value = 3001594.2
$cul.Range("M113").Value = -4166329.4
$cul.Range("N113").Value = -3005934.2

$gsm = $wb.Worksheets.Item("GSM")
$gsm.Range("H80").Value = 2490
$gsm.Range("I80").Value = 2488.3333
$gsm.Range("K80").Value = 2488.3333
$gsm.Range("M80").Value = -1490.3333

$gsm.Range("H83").Value = 2490
$gsm.Range("I83").Value = 2488.3333
$gsm.Range("K83").Value = 12441.6665
$gsm.Range("M83").Value = -7449.666499999999

$gsm.Range("H122").Value = 26298638
$gsm.Range("I122").Value = 35495850
$gsm.Range("J122").Value = 12502824
$gsm.Range("K122").Value = 106487550
$gsm.Range("L122").Value = 37508472
$gsm.Range("M122").Value = -106485100
$gsm.Range("N122").Value = -37513372

$gsm.Range("H123").Value = 19348.549
$gsm.Range("J123").Value = 19693.5
$gsm.Range("L123").Value = 19693.5
$gsm.Range("N123").Value = -24593.5

$gsm.Range("H132").Value = 2713.5715
$gsm.Range("I132").Value = 2790.6296
$gsm.Range("J132").Value = 2619
$gsm.Range("K132").Value = 8371.8888
$gsm.Range("L132").Value = 7857
$gsm.Range("M132").Value = -5841.888800000001
$gsm.Range("N132").Value = -12917

$ltw = $wb.Worksheets.Item("LTW")
$ltw.Range("H82").Value = 12376217
$ltw.Range("I82").Value = 2001278
$ltw.Range("J82").Value = 21022000
$ltw.Range("K82").Value = 2001278
$ltw.Range("L82").Value = 21022000
$ltw.Range("M82").Value = -2000917
$ltw.Range("N82").Value = -21022722

$ltw.Range("H85").Value = 12376217
$ltw.Range("I85").Value = 2001278
$ltw.Range("J85").Value = 21022000
$ltw.Range("K85").Value = 2001278
$ltw.Range("L85").Value = 21022000
$ltw.Range("M85").Value = -2000030
$ltw.Range("N85").Value = -21024496

$ltw.Range("H132").Value = 16669885
$ltw.Range("I132").Value = 30306500
$ltw.Range("J132").Value = 2911
$ltw.Range("K132").Value = 90919500
$ltw.Range("L132").Value = 8733
$ltw.Range("M132").Value = -90916970
$ltw.Range("N132").Value = -13793

$wvr = $wb.Worksheets.Item("WVR")
$wvr.Range("H81").Value = 1031.625
$wvr.Range("I81").Value = 1028
$wvr.Range("J81").Value = 1047.3334
$wvr.Range("K81").Value = 2056
$wvr.Range("L81").Value = 2094.6668
$wvr.Range("M81").Value = -995
$wvr.Range("N81").Value = -4216.6668

$wvr.Range("H84").Value = 1031.625
$wvr.Range("I84").Value = 1028
$wvr.Range("J84").Value = 1047.3334
$wvr.Range("K84").Value = 10280
$wvr.Range("L84").Value = 10473.334
$wvr.Range("M84").Value = -4976
$wvr.Range("N84").Value = -21081.334

$wvr.Range("H113").Value = 1194.925
$wvr.Range("J113").Value = 1292.0714
$wvr.Range("L113").Value = 3876.2142
$wvr.Range("N113").Value = -8216.2142
